# "added 4wk low sales check"
# Refresh the forecast numbers on "Forecast Comparison" (MyForecast, Inventory
# Coverage, Stockout Risk, Seasonality Index) and roll the new totals up into
# the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H),
#     Stockout Risk (I), Seasonality Index (L) ---------------------------

# Week W10 (row 2)
$ws1.Range("D2").Value = 58
$ws1.Range("H2").Value = 7.17
$ws1.Range("L2").Value = 0.94

# Week W11 (row 3)
$ws1.Range("D3").Value = 56
$ws1.Range("H3").Value = 6.46
$ws1.Range("L3").Value = 1.01

# Week W12 (row 4)
$ws1.Range("D4").Value = 58
$ws1.Range("H4").Value = 5.22
$ws1.Range("L4").Value = 1.19

# Week W13 (row 5)
$ws1.Range("D5").Value = 77
$ws1.Range("H5").Value = 3.19
$ws1.Range("L5").Value = 1.04

# Week W14 (row 6)
$ws1.Range("D6").Value = 94
$ws1.Range("H6").Value = 1.79
$ws1.Range("L6").Value = 1.2

# Week W15 (row 7)
$ws1.Range("D7").Value = 98
$ws1.Range("H7").Value = 0.76
$ws1.Range("I7").Value = "Low"
$ws1.Range("L7").Value = 0.93

# Week W16 (row 8)
$ws1.Range("D8").Value = 88
$ws1.Range("L8").Value = 1.17

# Week W17 (row 9)
$ws1.Range("D9").Value = 78
$ws1.Range("L9").Value = 0.86

# Week W18 (row 10)
$ws1.Range("D10").Value = 77
$ws1.Range("L10").Value = 0.94

# Week W19 (row 11)
$ws1.Range("D11").Value = 83
$ws1.Range("L11").Value = 1.06

# Week W20 (row 12)
$ws1.Range("D12").Value = 88
$ws1.Range("L12").Value = 0.94

# Week W21 (row 13)
$ws1.Range("D13").Value = 86
$ws1.Range("L13").Value = 0.86

# Week W22 (row 14)
$ws1.Range("D14").Value = 81
$ws1.Range("L14").Value = 0.81

# Week W23 (row 15)
$ws1.Range("D15").Value = 79
$ws1.Range("L15").Value = 1.03

# Week W24 (row 16)
$ws1.Range("D16").Value = 80
$ws1.Range("L16").Value = 1.08

# Week W25 (row 17)
$ws1.Range("D17").Value = 78
$ws1.Range("L17").Value = 1.06

# --- Summary: roll the refreshed forecast back into the totals -----------
# (values are stored as text on this sheet, so keep them text with a
# leading apostrophe, matching how the sheet was originally authored)

$ws2.Range("B9").Value  = "'1267"   # Total Forecast (16 Weeks)
$ws2.Range("B10").Value = "'612"    # Total Forecast (8 Weeks)
$ws2.Range("B11").Value = "'251"    # Total Forecast (4 Weeks)
$ws2.Range("B12").Value = "'98"     # Max Forecast
$ws2.Range("B14").Value = "'56"     # Min Forecast
